$wb = $excel.ActiveWorkbook
$cleaned = $wb.Worksheets.Item("cleaned")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $cleaned)
$newSheet.Name = "Sheet1"
Write-Output ($newSheet | Get-Member -Name "*Height*" | Out-String)
